# "updated task used in testing"
# Row 4 of the training schedule is revised: the step-count columns (D, F)
# move from 2 -> 3 and the derived/packed value in H moves from 36 -> 46.
# The user's active selection also moves up a row, from D5 to D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 4).Value = 3   # D4: x_corrSteps  2 -> 3
$ws.Cells.Item(4, 6).Value = 3   # F4: y_corrSteps  2 -> 3
$ws.Cells.Item(4, 8).Value = 46  # H4: alienID     36 -> 46

# Move the selection/active cell to D4 (was D5).
$ws.Range("D4").Select()
